# Scheduled market-data refresh for Moogle Treasure Trove leve profit sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) per row
# with freshly pulled Universalis price data; some rows gain/lose M/N cells
# when HQ pricing data becomes available/unavailable for that item.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 627.0952
$ws.Range("I9").Value = 350
$ws.Range("K9").Value = 350
$ws.Range("M9").Value = -181

# Row 33
$ws.Range("H33").Value = 1114.2354
$ws.Range("I33").Value = 616.2857
$ws.Range("J33").Value = 3438
$ws.Range("K33").Value = 616.2857
$ws.Range("L33").Value = 3438
$ws.Range("M33").Value = -387.2857
$ws.Range("N33").Value = -3896

# Row 112
$ws.Range("H112").Value = 11529.2
$ws.Range("J112").Value = 12589
$ws.Range("L112").Value = 37767
$ws.Range("N112").Value = -39983

# Row 127
$ws.Range("H127").Value = 2128.3333
$ws.Range("I127").Value = 921.625
$ws.Range("K127").Value = 2764.875
$ws.Range("M127").Value = 2195.125

# Row 131
$ws.Range("H131").Value = 3247.3333
$ws.Range("I131").Value = 3096.8
$ws.Range("K131").Value = 9290.400000000001
$ws.Range("M131").Value = -4250.400000000001

# Row 137
$ws.Range("H137").Value = 3479.561
$ws.Range("I137").Value = 3067.543
$ws.Range("J137").Value = 5883
$ws.Range("K137").Value = 9202.629000000001
$ws.Range("L137").Value = 17649
$ws.Range("M137").Value = -6652.629000000001
$ws.Range("N137").Value = -22749

# Row 138
$ws.Range("H138").Value = 5370.155
$ws.Range("I138").Value = 5599.136
$ws.Range("K138").Value = 16797.408
$ws.Range("M138").Value = -11657.408

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 83755.336
$ws.Range("J24").Value = 83755.336
$ws.Range("L24").Value = 83755.336
$ws.Range("N24").Value = -84503.336

# Row 32
$ws.Range("H32").Value = 12076.404
$ws.Range("I32").Value = 9182.519
$ws.Range("K32").Value = 9182.519
$ws.Range("M32").Value = -8895.519

# Row 43
$ws.Range("H43").Value = 17265.666
$ws.Range("I43").Value = 12000
$ws.Range("J43").Value = 19898.5
$ws.Range("K43").Value = 12000
$ws.Range("L43").Value = 19898.5
$ws.Range("M43").Value = -11687
$ws.Range("N43").Value = -20524.5

# Row 45
$ws.Range("H45").Value = 3492.8235
$ws.Range("I45").Value = 2743.75
$ws.Range("J45").Value = 4158.6665
$ws.Range("K45").Value = 2743.75
$ws.Range("L45").Value = 4158.6665
$ws.Range("M45").Value = -2366.75
$ws.Range("N45").Value = -4912.6665

# Row 97
$ws.Range("H97").Value = 70
$ws.Range("I97").Value = 70
$ws.Range("K97").Value = 70
$ws.Range("M97").Value = 426

# Row 100
$ws.Range("H100").Value = 83755.336
$ws.Range("J100").Value = 83755.336
$ws.Range("L100").Value = 83755.336
$ws.Range("N100").Value = -85919.336

$ws = $wb.Worksheets.Item("BSM")
# Row 100
$ws.Range("H100").Value = 44128.2
$ws.Range("J100").Value = 44128.2
$ws.Range("L100").Value = 44128.2
$ws.Range("N100").Value = -46292.2

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 4540.484
$ws.Range("I122").Value = 4340.7
$ws.Range("J122").Value = 4903.727
$ws.Range("K122").Value = 13022.1
$ws.Range("L122").Value = 14711.181
$ws.Range("M122").Value = -10572.1
$ws.Range("N122").Value = -19611.181

# Row 141
$ws.Range("H141").Value = 283015.22
$ws.Range("J141").Value = 293182.78
$ws.Range("L141").Value = 293182.78
$ws.Range("N141").Value = -303542.78

$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 12431.1875
$ws.Range("I87").Value = 4724.75
$ws.Range("K87").Value = 14174.25
$ws.Range("M87").Value = -12926.25

# Row 90
$ws.Range("H90").Value = 12431.1875
$ws.Range("I90").Value = 4724.75
$ws.Range("K90").Value = 43497
$ws.Range("M90").Value = -36282.75

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3675.1924
$ws.Range("I102").Value = 3407.6667
$ws.Range("J102").Value = 4798.8
$ws.Range("K102").Value = 3407.6667
$ws.Range("L102").Value = 4798.8
$ws.Range("M102").Value = -1785.6667
$ws.Range("N102").Value = -8042.8

# Row 122
$ws.Range("H122").Value = 8769.5625
$ws.Range("J122").Value = 13856.714
$ws.Range("L122").Value = 41570.142
$ws.Range("N122").Value = -46470.142

# Row 126
$ws.Range("H126").Value = 8105.3213
$ws.Range("J126").Value = 8102.5
$ws.Range("L126").Value = 24307.5
$ws.Range("N126").Value = -29247.5

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5453.7896
$ws.Range("I40").Value = 4101.375
$ws.Range("K40").Value = 4101.375
$ws.Range("M40").Value = -3965.375

# Row 46
$ws.Range("H46").Value = 3622.15
$ws.Range("I46").Value = 3869.125
$ws.Range("J46").Value = 3457.5
$ws.Range("K46").Value = 3869.125
$ws.Range("L46").Value = 3457.5
$ws.Range("M46").Value = -3681.125
$ws.Range("N46").Value = -3833.5

# Row 55
$ws.Range("H55").Value = 596.73334
$ws.Range("I55").Value = 301.9
$ws.Range("J55").Value = 1186.4
$ws.Range("K55").Value = 301.9
$ws.Range("L55").Value = 1186.4
$ws.Range("M55").Value = -128.9
$ws.Range("N55").Value = -1532.4

# Row 68
$ws.Range("H68").Value = 6977
$ws.Range("I68").Value = 6977
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 6977
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -6228
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 6977
$ws.Range("I71").Value = 6977
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 34885
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -31141
$ws.Range("N71").ClearContents()

# Row 132
$ws.Range("H132").Value = 4051.5
$ws.Range("I132").Value = 3366.1936
$ws.Range("K132").Value = 10098.5808
$ws.Range("M132").Value = -7568.5808

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 20583
$ws.Range("J15").Value = 20475
$ws.Range("L15").Value = 20475
$ws.Range("N15").Value = -21051

# Row 18
$ws.Range("H18").Value = 200000
$ws.Range("I18").Value = 200000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 200000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -199827
$ws.Range("N18").ClearContents()

# Row 80
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996

# Row 83
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984

# Row 96
$ws.Range("H96").Value = 23451.5
$ws.Range("J96").Value = 45000
$ws.Range("L96").Value = 45000
$ws.Range("N96").Value = -47746

# Row 126
$ws.Range("H126").Value = 3361
$ws.Range("I126").Value = 2368.2307
$ws.Range("K126").Value = 7104.6921
$ws.Range("M126").Value = -4634.6921
